$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 138.68965
$ws.Range("I2").Value = 98.61905
$ws.Range("J2").Value = 243.875
$ws.Range("K2").Value = 98.61905
$ws.Range("L2").Value = 243.875
$ws.Range("M2").Value = 14.38095
$ws.Range("N2").Value = -469.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1353.2941
$ws.Range("I137").Value = 1144.48
$ws.Range("K137").Value = 3433.44
$ws.Range("M137").Value = -883.4400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20444.543
$ws.Range("I32").Value = 22299.627
$ws.Range("K32").Value = 22299.627
$ws.Range("M32").Value = -22012.627

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2961.6924
$ws.Range("I61").Value = 2562.75
$ws.Range("K61").Value = 2562.75
$ws.Range("M61").Value = -2350.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 38000
$ws.Range("J82").Value = 38000
$ws.Range("L82").Value = 38000
$ws.Range("N82").Value = -38722

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 38000
$ws.Range("J85").Value = 38000
$ws.Range("L85").Value = 38000
$ws.Range("N85").Value = -40496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 19800
$ws.Range("J86").Value = 19800
$ws.Range("L86").Value = 19800
$ws.Range("N86").Value = -22172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H87").Value = 38000
$ws.Range("J87").Value = 38000
$ws.Range("L87").Value = 38000
$ws.Range("N87").Value = -40496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 19800
$ws.Range("J89").Value = 19800
$ws.Range("L89").Value = 59400
$ws.Range("N89").Value = -71256

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H90").Value = 38000
$ws.Range("J90").Value = 38000
$ws.Range("L90").Value = 114000
$ws.Range("N90").Value = -126480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 953.9167
$ws.Range("I97").Value = 858.8182
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 858.8182
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -362.8182
$ws.Range("N97").Value = -2992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3265.3635
$ws.Range("I122").Value = 3331
$ws.Range("J122").Value = 3227.8572
$ws.Range("K122").Value = 9993
$ws.Range("L122").Value = 9683.571599999999
$ws.Range("M122").Value = -7543
$ws.Range("N122").Value = -14583.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1623.9032
$ws.Range("I132").Value = 1185.7727
$ws.Range("J132").Value = 2694.889
$ws.Range("K132").Value = 3557.3181
$ws.Range("L132").Value = 8084.667
$ws.Range("M132").Value = -1027.3181
$ws.Range("N132").Value = -13144.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2961.6924
$ws.Range("I136").Value = 2562.75
$ws.Range("K136").Value = 7688.25
$ws.Range("M136").Value = -5138.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2603.0303
$ws.Range("I132").Value = 1971.6666
$ws.Range("J132").Value = 3707.9167
$ws.Range("K132").Value = 5914.9998
$ws.Range("L132").Value = 11123.7501
$ws.Range("M132").Value = -3384.9998
$ws.Range("N132").Value = -16183.7501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 949.9487
$ws.Range("I134").Value = 900.6667
$ws.Range("J134").Value = 1221
$ws.Range("K134").Value = 2702.0001
$ws.Range("L134").Value = 3663
$ws.Range("M134").Value = -167.0001000000002
$ws.Range("N134").Value = -8733

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 643.913
$ws.Range("I5").Value = 454.16666
$ws.Range("K5").Value = 1362.49998
$ws.Range("M5").Value = -1250.49998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 9000
$ws.Range("N80").Value = -10872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 27000
$ws.Range("N83").Value = -36360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1512.1177
$ws.Range("I134").Value = 1362
$ws.Range("K134").Value = 4086
$ws.Range("M134").Value = 984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 643.913
$ws.Range("I135").Value = 454.16666
$ws.Range("K135").Value = 4087.49994
$ws.Range("M135").Value = -1552.49994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2571.3572
$ws.Range("I122").Value = 1624.75
$ws.Range("K122").Value = 4874.25
$ws.Range("M122").Value = -2424.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5561318.5
$ws.Range("I126").Value = 6916
$ws.Range("J126").Value = 33333332
$ws.Range("K126").Value = 20748
$ws.Range("L126").Value = 99999996
$ws.Range("M126").Value = -18278
$ws.Range("N126").Value = -100004936

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1293.5294
$ws.Range("I61").Value = 1092.1428
$ws.Range("J61").Value = 2233.3333
$ws.Range("K61").Value = 1092.1428
$ws.Range("L61").Value = 2233.3333
$ws.Range("M61").Value = -890.1428000000001
$ws.Range("N61").Value = -2637.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1293.5294
$ws.Range("I113").Value = 1092.1428
$ws.Range("J113").Value = 2233.3333
$ws.Range("K113").Value = 1092.1428
$ws.Range("L113").Value = 2233.3333
$ws.Range("M113").Value = 1077.8572
$ws.Range("N113").Value = -6573.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 15785.4
$ws.Range("I132").Value = 20964.295
$ws.Range("J132").Value = 4780.25
$ws.Range("K132").Value = 62892.88499999999
$ws.Range("L132").Value = 14340.75
$ws.Range("M132").Value = -60362.88499999999
$ws.Range("N132").Value = -19400.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7601.1763
$ws.Range("I136").Value = 9016.923000000001
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 27050.769
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -24500.769
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 568
$ws.Range("J113").Value = 450
$ws.Range("L113").Value = 1350
$ws.Range("N113").Value = -5690

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2132.389
$ws.Range("J122").Value = 1459.375
$ws.Range("L122").Value = 4378.125
$ws.Range("N122").Value = -9278.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1538.8
$ws.Range("I126").Value = 1538.8
$ws.Range("K126").Value = 4616.4
$ws.Range("M126").Value = -2146.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1749.0416
$ws.Range("I132").Value = 1209.7222
$ws.Range("J132").Value = 3367
$ws.Range("K132").Value = 3629.1666
$ws.Range("L132").Value = 10101
$ws.Range("M132").Value = -1099.1666
$ws.Range("N132").Value = -15161
